$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.538.49'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '2.623.34'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.51'
$ws.Range("E5").Value = '  +3.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.81'
$ws.Range("E6").Value = '  +2.26%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +0.83%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("E10").Value = '  +1.51%  '
$ws.Range("E11").Value = '  +2.48%  '
$ws.Range("E12").Value = '  +1.90%  '
$ws.Range("D13").Value = '3.084.55'
$ws.Range("E13").Value = '  +1.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.34'
$ws.Range("E14").Value = '  +13.47%  '
$ws.Range("D15").Value = '60.529.69'
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("E16").Value = '  +1.78%  '
$ws.Range("D17").Value = '2.628.14'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.52'
$ws.Range("E18").Value = '  +2.62%  '
$ws.Range("E19").Value = '  +2.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '348.88'
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("E21").Value = '  -0.55%  '
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.528'
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("E24").Value = '  +1.56%  '
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("E26").Value = '  +1.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.17'
$ws.Range("E27").Value = '  +7.63%  '
$ws.Range("E28").Value = '  +16.33%  '
$ws.Range("E29").Value = '  +3.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.82'
$ws.Range("E31").Value = '  +5.59%  '
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.45'
$ws.Range("E34").Value = '  +6.02%  '
$ws.Range("E35").Value = '  +9.53%  '
$ws.Range("E36").Value = '  +9.06%  '
$ws.Range("E37").Value = '  +4.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '330.85'
$ws.Range("E38").Value = '  +13.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.89'
$ws.Range("E39").Value = '  +3.51%  '
$ws.Range("E40").Value = '  +5.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.871'
$ws.Range("E41").Value = '  +2.05%  '
$ws.Range("E42").Value = '  +8.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '133.60'
$ws.Range("E43").Value = '  -2.73%  '
$ws.Range("E44").Value = '  +2.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.12'
$ws.Range("E45").Value = '  +3.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.25%  '
$ws.Range("E47").Value = '  +2.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.611'
$ws.Range("E48").Value = '  +1.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.36'
$ws.Range("E49").Value = '  +4.41%  '
$ws.Range("E50").Value = '  +2.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.73'
$ws.Range("E51").Value = '  +0.68%  '
